# Generate Report for Handoff
# Regenerated the handoff XLIFF files for the four records that were still
# pending ("low" priority, stale handoff timestamp) so they now carry the
# "ht" priority and a refreshed handoff timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 4-7 correspond to: 20d4192d, 59318f4a, 9837d349, ff75e68d
$rows = 4,5,6,7

foreach ($r in $rows) {
    # Priority moved from "low" to "ht" on both locale sheets.
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # zh-cn's "Latest Handoff Datetime" got its own, independent refresh.
    $zhcn.Range("H$r").Value = "2016-09-05 20:37:35"

    # de-de's "Latest Handoff Datetime" shares the same underlying
    # timestamp as the Overview sheet's "Latest HO Xliff Generate Date".
    $dede.Range("H$r").Value = "2016-09-05 20:37:40"

    $overview.Range("G$r").Value = "2016-09-05 20:37:40"
}
